# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.445.89"
$ws.Range("E2").Value = "  -2.50%  "
$ws.Range("D3").Value = "2.337.33"
$ws.Range("E3").Value = "  -3.28%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'318.51"
$ws.Range("E5").Value = "  -1.14%  "
$ws.Range("D6").Value = "'104.36"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  -1.54%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -6.41%  "
$ws.Range("D10").Value = "'40.46"
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("E11").Value = "  -2.85%  "
$ws.Range("D12").Value = "'8.29"
$ws.Range("E12").Value = "  -4.64%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.983"
$ws.Range("E13").Value = "  -5.05%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "'0.105"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").Value = "'15.80"
$ws.Range("E15").Value = "  -9.00%  "
$ws.Range("D16").Value = "2.690.72"
$ws.Range("E16").Value = "  -3.18%  "
$ws.Range("D17").Value = "2.326.95"
$ws.Range("E17").Value = "  -6.83%  "
$ws.Range("D18").Value = "42.452.34"
$ws.Range("E18").Value = "  -2.56%  "
$ws.Range("D19").Value = "'7.69"
$ws.Range("E19").Value = "  +3.83%  "
$ws.Range("E20").Value = "  -4.23%  "
$ws.Range("D21").Value = "'76.85"
$ws.Range("E21").Value = "  +1.73%  "
$ws.Range("D22").Value = "'3.54"
$ws.Range("E22").Value = "  +2.24%  "
$ws.Range("D23").Value = "'259.18"
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("E24").Value = "  -5.32%  "
$ws.Range("D25").Value = "'9.61"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'11.34"
$ws.Range("E27").Value = "  -5.47%  "
$ws.Range("D28").Value = "'23.00"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").Value = "'174.62"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").Value = "'35.69"
$ws.Range("E31").Value = "  -6.49%  "
$ws.Range("D32").Value = "'0.0888"
$ws.Range("E32").Value = "  -5.19%  "
$ws.Range("E33").Value = "  -7.88%  "
$ws.Range("D34").Value = "'6.03"
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("E36").Value = "  +5.69%  "
$ws.Range("D37").Value = "'4.51"
$ws.Range("E37").Value = "  -7.53%  "
$ws.Range("D38").Value = "'0.0354"
$ws.Range("E38").Value = "  -4.56%  "
$ws.Range("E39").Value = "  -5.35%  "
$ws.Range("D40").Value = "'2.62"
$ws.Range("E40").Value = "  -9.86%  "
$ws.Range("E41").Value = "  -11.19%  "
$ws.Range("D42").Value = "'70.44"
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("D43").Value = "'0.232"
$ws.Range("E43").Value = "  -1.12%  "
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "'115.28"
$ws.Range("E45").Value = "  -7.63%  "
$ws.Range("D46").Value = "'11.74"
$ws.Range("E46").Value = "  -7.04%  "
$ws.Range("D47").Value = "'5.48"
$ws.Range("E47").Value = "  -3.63%  "
$ws.Range("D48").Value = "'9.13"
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("D49").Value = "'84.71"
$ws.Range("E49").Value = "  +9.78%  "
$ws.Range("D50").Value = "'72.94"
$ws.Range("E50").Value = "  +2.60%  "
$ws.Range("D51").Value = "'0.0996"
$ws.Range("E51").Value = "  -1.34%  "
